$wb = $excel.ActiveWorkbook

# Update "想去人数" (F column) values on both the "展览" sheet and the
# "全部类型" sheet, which contain duplicated data.
$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F2").Value = 100
    $ws.Range("F3").Value = 21
}
